$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string label updates (country reorder) ---
$ws.Range("A97").Value = "Malasia"
$ws.Range("A98").Value = "Namibia"
$ws.Range("A159").Value = "Letonia"
$ws.Range("A160").Value = "Togo"
$ws.Range("A161").Value = "Republica de Chipre"
$ws.Range("A215").Value = "Montserrat"
$ws.Range("A216").Value = "Islas Malvinas"

# --- Timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Septiembre de 2020 a las 13:32"

# --- Numeric data updates (updated COVID case counts) ---
# Row 16
$ws.Range("B16").Value = 457219
$ws.Range("C16").Value = 3582
$ws.Range("D16").Value = 380956
$ws.Range("E16").Value = 50094
$ws.Range("G16").Value = 183
$ws.Range("H16").Value = 26169

# Row 18
$ws.Range("B18").Value = 363479
$ws.Range("C18").Value = 1436
$ws.Range("D18").Value = 275487
$ws.Range("E18").Value = 82741
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = 5251

# Row 33
$ws.Range("B33").Value = 125760
$ws.Range("C33").Value = 227
$ws.Range("D33").Value = 122699
$ws.Range("E33").Value = 2847

# Row 50
$ws.Range("B50").Value = 77817
$ws.Range("C50").Value = 1559
$ws.Range("D50").Value = 56428
$ws.Range("E50").Value = 20891
$ws.Range("G50").Value = 7
$ws.Range("H50").Value = 498

# Row 61
$ws.Range("B61").Value = 53282
$ws.Range("C61").Value = 411
$ws.Range("E61").Value = 8508
$ws.Range("G61").Value = 5
$ws.Range("H61").Value = 2074

# Row 90
$ws.Range("B90").Value = 16408
$ws.Range("C90").Value = 31
$ws.Range("D90").Value = 15301
$ws.Range("E90").Value = 877
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 230

# Row 91
$ws.Range("B91").Value = 14982
$ws.Range("C91").Value = 37
$ws.Range("D91").Value = 12437
$ws.Range("E91").Value = 2234
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 311

# Row 97
$ws.Range("B97").Value = 11224
$ws.Range("C97").Value = 89
$ws.Range("D97").Value = 9967
$ws.Range("E97").Value = 1121
$ws.Range("G97").Value = 2
$ws.Range("H97").Value = 136

# Row 98
$ws.Range("B98").Value = 11140
$ws.Range("D98").Value = 8937
$ws.Range("E98").Value = 2082
$ws.Range("H98").Value = 121

# Row 116
$ws.Range("E116").Value = 3033
$ws.Range("G116").Value = 3
$ws.Range("H116").Value = 39

# Row 146
$ws.Range("B146").Value = 3058
$ws.Range("C146").Value = 23
$ws.Range("D146").Value = 2562
$ws.Range("E146").Value = 462

# Row 159
$ws.Range("B159").Value = 1824
$ws.Range("C159").Value = 95
$ws.Range("D159").Value = 1307
$ws.Range("E159").Value = 480
$ws.Range("H159").Value = 37

# Row 160
$ws.Range("B160").Value = 1759
$ws.Range("D160").Value = 1341
$ws.Range("E160").Value = 370
$ws.Range("H160").Value = 48

# Row 161
$ws.Range("B161").Value = 1743
$ws.Range("D161").Value = 1369
$ws.Range("E161").Value = 352
$ws.Range("H161").Value = 22

# Row 168
$ws.Range("D168").Value = 1010
$ws.Range("E168").Value = 49

# Row 215
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

# Row 216
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
